$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 72 does not exist yet; create it with the same formatting as row 71
$ws.Range("A71:C71").Copy()
$ws.Range("A72:C72").PasteSpecial(-4122)

# Row 3
$ws.Range("C3").Value = 1050.0

# Row 4
$ws.Range("C4").Value = 731.0

# Row 5
$ws.Range("C5").Value = 461.0

# Row 9
$ws.Range("A9").Value = 'Regione Lombardia'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '80050050154'
$ws.Range("C9").Value = 368.0

# Row 10
$ws.Range("A10").Value = 'Maggioli SPA'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '06188330150'
$ws.Range("C10").Value = 360.0

# Row 11
$ws.Range("A11").Value = 'Provincia Autonoma di Trento'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = '00337460224'
$ws.Range("C11").Value = 331.0

# Row 12
$ws.Range("A12").Value = 'Banca Popolare di Sondrio, Società Cooperativa per Azioni'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = '00053810149'
$ws.Range("C12").Value = 304.0

# Row 13
$ws.Range("C13").Value = 230.0

# Row 14
$ws.Range("A14").Value = 'Alto Adige Riscossioni Spa'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = '02805390214'
$ws.Range("C14").Value = 180.0

# Row 15
$ws.Range("A15").Value = 'Progetti e Soluzioni SPA'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = '06423240727'
$ws.Range("C15").Value = 164.0

# Row 16
$ws.Range("A16").Value = 'P.A. Digitale spa'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = '06628860964'
$ws.Range("C16").Value = 143.0

# Row 17
$ws.Range("A17").Value = 'APKAPPA S.R.L.'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = '08543640158'
$ws.Range("C17").Value = 141.0

# Row 18
$ws.Range("A18").Value = 'Regione Marche'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = '80008630420'
$ws.Range("C18").Value = 139.0

# Row 19
$ws.Range("A19").Value = 'ADVANCED SYSTEMS srl'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = '03383350638'
$ws.Range("C19").Value = 123.0

# Row 20
$ws.Range("C20").Value = 117.0

# Row 21
$ws.Range("A21").Value = 'Regione Basilicata'
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = '80002950766'
$ws.Range("C21").Value = 106.0

# Row 22
$ws.Range("A22").Value = 'Regione Toscana'
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '01386030488'
$ws.Range("C22").Value = 100.0

# Row 23
$ws.Range("A23").Value = 'ANCITEL'
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '07196850585'
$ws.Range("C23").Value = 90.0

# Row 24
$ws.Range("A24").Value = 'Next Step Solution'
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = '02554480349'
$ws.Range("C24").Value = 74.0

# Row 25
$ws.Range("A25").Value = 'Regione Autonoma Friuli-Venezia Giulia'
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '80014930327'
$ws.Range("C25").Value = 73.0

# Row 26
$ws.Range("A26").Value = 'Regione Piemonte'
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = '80087670016'
$ws.Range("C26").Value = 63.0

# Row 27
$ws.Range("C27").Value = 59.0

# Row 32
$ws.Range("C32").Value = 43.0

# Row 35
$ws.Range("A35").Value = 'UNIMATICA S.P.A'
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = '02098391200'
$ws.Range("C35").Value = 38.0

# Row 36
$ws.Range("A36").Value = 'SI.net Servizi Informatici S.r.L.'
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = '02743730125'
$ws.Range("C36").Value = 35.0

# Row 37
$ws.Range("A37").Value = 'Unicredit, Societa'' per Azioni'
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = '00348170101'
$ws.Range("C37").Value = 33.0

# Row 38
$ws.Range("A38").Value = 'ROMA CAPITALE'
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = '02438750586'
$ws.Range("C38").Value = 32.0

# Row 39
$ws.Range("A39").Value = 'PMPay s.r.l.'
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = '08747230962'
$ws.Range("C39").Value = 29.0

# Row 40
$ws.Range("A40").Value = 'DCS SOFTWARE E SERVIZI S.R.L.'
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = '08063140019'
$ws.Range("C40").Value = 24.0

# Row 42
$ws.Range("A42").Value = 'Siscom SPA'
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = '01778000040'
$ws.Range("C42").Value = 19.0

# Row 43
$ws.Range("A43").Value = 'Regione Lazio'
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = '80143490581'

# Row 44
$ws.Range("A44").Value = 'Citta'' Metropolitana di Roma Capitale'
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = '80034390585'

# Row 45
$ws.Range("A45").Value = 'Nexi SpA'
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = '13212880150'
$ws.Range("C45").Value = 18.0

# Row 47
$ws.Range("C47").Value = 17.0

# Row 48
$ws.Range("A48").Value = 'Novares Spa'
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = '12105121003'

# Row 49
$ws.Range("A49").Value = 'Servizi Locali SpA'
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = '03170580751'

# Row 51
$ws.Range("A51").Value = 'Crédit Agricole Group Solutions Società Consortile per azioni'
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = '02771790348'
$ws.Range("C51").Value = 10.0

# Row 52
$ws.Range("A52").Value = 'Aric Agenzia Regionale di Informatica e Committenza'
$ws.Range("B52").NumberFormat = "@"
$ws.Range("B52").Value = '91022630676'

# Row 53
$ws.Range("A53").Value = 'UBI Banca'
$ws.Range("B53").NumberFormat = "@"
$ws.Range("B53").Value = '03053920165'

# Row 55
$ws.Range("A55").Value = 'Numera Sistemi e Informatica SpA'
$ws.Range("B55").NumberFormat = "@"
$ws.Range("B55").Value = '01265230902'
$ws.Range("C55").Value = 6.0

# Row 57
$ws.Range("A57").Value = 'Comune di Catania'
$ws.Range("B57").NumberFormat = "@"
$ws.Range("B57").Value = '00137020871'

# Row 58
$ws.Range("A58").Value = 'ISWEB S.p.A.'
$ws.Range("B58").NumberFormat = "@"
$ws.Range("B58").Value = '01722270665'
$ws.Range("C58").Value = 3.0

# Row 60
$ws.Range("A60").Value = 'Engineering Ingegneria Informatica SpA'
$ws.Range("B60").NumberFormat = "@"
$ws.Range("B60").Value = '00967720285'
$ws.Range("C60").Value = 1.0

# Row 61
$ws.Range("A61").Value = 'Agenzia Italiana del Farmaco - AIFA'
$ws.Range("B61").NumberFormat = "@"
$ws.Range("B61").Value = '97345810580'

# Row 62
$ws.Range("A62").Value = 'Banco BPM Società per Azioni'
$ws.Range("B62").NumberFormat = "@"
$ws.Range("B62").Value = '09722490969'

# Row 63
$ws.Range("A63").Value = 'ARGO SOFTWARE SRL'
$ws.Range("B63").NumberFormat = "@"
$ws.Range("B63").Value = '00838520880'

# Row 64
$ws.Range("A64").Value = 'MegASP S.r.l.'
$ws.Range("B64").NumberFormat = "@"
$ws.Range("B64").Value = '09898030151'

# Row 65
$ws.Range("A65").Value = 'Società Almaviva S.p.A.'
$ws.Range("B65").NumberFormat = "@"
$ws.Range("B65").Value = '08450891000'

# Row 66
$ws.Range("A66").Value = 'San Marco SPA'
$ws.Range("B66").NumberFormat = "@"
$ws.Range("B66").Value = '04142440728'

# Row 67
$ws.Range("A67").Value = 'I.C.A. - Imposte Comunali Affini – s.r.l.'
$ws.Range("B67").NumberFormat = "@"
$ws.Range("B67").Value = '02478610583'

# Row 68
$ws.Range("A68").Value = 'ICCREA Banca SpA'
$ws.Range("B68").NumberFormat = "@"
$ws.Range("B68").Value = '04774801007'

# Row 69
$ws.Range("A69").Value = 'ARCA Servizi s.r.l'
$ws.Range("B69").NumberFormat = "@"
$ws.Range("B69").Value = '09106071005'

# Row 70
$ws.Range("A70").Value = 'Ministero dello Sviluppo Economico'
$ws.Range("B70").NumberFormat = "@"
$ws.Range("B70").Value = '80230390587'

# Row 71
$ws.Range("A71").Value = 'Softline srl'
$ws.Range("B71").NumberFormat = "@"
$ws.Range("B71").Value = '12299030150'

# Row 72
$ws.Range("A72").Value = 'CityPoste Payment Digital S.r.l.'
$ws.Range("B72").NumberFormat = "@"
$ws.Range("B72").Value = '02003750672'
$ws.Range("C72").Value = 1.0
